# Insert a new data row at row 169 (shifting existing rows 169-295 down to
# 170-296) and populate it with the new record described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 45072
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100114007
$ws.Cells.Item(169, 7).Value = "Jengibre"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 100
$ws.Cells.Item(169, 11).Value = 24000
$ws.Cells.Item(169, 12).Value = 24000
$ws.Cells.Item(169, 13).Value = 24000
$ws.Cells.Item(169, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(169, 15).Value = "Perú"
$ws.Cells.Item(169, 16).Value = 1846
$ws.Cells.Item(169, 17).Value = 13
$ws.Cells.Item(169, 18).Value = "Hortaliza"
